# GPLIM-1: Remove empty column that was supposed to be moved.
#
# Column J ("J:J") on both worksheets (P-RNA-0004 and P-EXEX-0001) is
# completely empty -- it was left behind when data that used to live there
# was moved elsewhere. Delete it so the columns to its right (Sort Column,
# the product/billing quantity columns, Billing Errors, etc.) shift one
# position to the left, closing the gap.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("P-RNA-0004")
$ws2 = $wb.Worksheets.Item("P-EXEX-0001")

# Remove the stray empty column on each sheet.
$ws1.Range("J1").EntireColumn.Delete() | Out-Null
$ws2.Range("J1").EntireColumn.Delete() | Out-Null

# Restore the selection on each sheet to the (now-shifted) column J, and
# leave P-EXEX-0001 as the active/selected tab.
$ws1.Range("J1:J1048576").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("J1:J1048576").Select() | Out-Null
